$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert an extra row to make room for a new explanatory bullet point ---
# (old rows 14 and 15 shift down to 15 and 16)
$null = $ws.Rows.Item(14).Insert()

# --- Rewrite the "notes" bullet list (rows 12-16) ---
$ws.Range("B13").Value = "* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message."
$ws.Range("B14").Value = "    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation."
$ws.Range("B16").Value = "* You cannot change picture using ""ope=changepic"". This is because drawing information are not saved directly in the sheet."

# --- Add a new "Score" column to the Example #1 table ---
# Header cell (row 20), formatted like the other header cells (D20)
$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "Score"

# TBS field cell (row 21), formatted like the other field cells (D21)
$ws.Range("D21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E21").Value = "[a.score;ope=xlsxNum]"

# --- Add a "Total:" row (row 19) that sums the Score column ---
$ws.Range("D19").Value = "Total:"
$ws.Range("D19").HorizontalAlignment = -4152
$ws.Range("E19").Formula = "=SUM(E21:E2000)"
$ws.Range("E19").Font.Bold = $true
$ws.Range("E19").NumberFormat = "#,##0.0"

# Finish formatting the new Score data cell: number format + right alignment
$ws.Range("E21").NumberFormat = "#,##0.0"
$ws.Range("E21").HorizontalAlignment = -4152

$null = $ws.Range("E20").Select()

Write-Host "done"
